$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "530.54", "57.003.31") are preserved exactly as text, matching the source data.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "57.003.31"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "2.320.25"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").Value = "530.54"
$ws.Range("E5").Value = "  +1.84%  "

$ws.Range("D6").Value = "132.48"

$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").Value = "2.347.47"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("E12").Value = "  -2.70%  "

$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("D14").Value = "2.741.42"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").Value = "23.47"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("D16").Value = "57.046.47"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("D18").Value = "2.338.66"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").Value = "337.26"
$ws.Range("E19").Value = "  +2.38%  "

$ws.Range("D20").Value = "10.43"
$ws.Range("E20").Value = "  -1.69%  "

$ws.Range("D21").Value = "6.89"
$ws.Range("E21").Value = "  +2.38%  "

$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -1.90%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "61.61"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").Value = "8.72"
$ws.Range("E26").Value = "  +1.15%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("E28").Value = "  +1.26%  "

$ws.Range("D29").Value = "173.76"
$ws.Range("E29").Value = "  +3.40%  "

$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("E31").Value = "  -2.56%  "

$ws.Range("E32").Value = "  -2.92%  "

$ws.Range("D33").Value = "18.51"
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  -3.65%  "

$ws.Range("D37").Value = "0.916"
$ws.Range("E37").Value = "  -0.89%  "

$ws.Range("D38").Value = "3.99"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("D39").Value = "39.28"
$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("E40").Value = "  -2.32%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "149.07"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.61"
$ws.Range("E42").Value = "  +4.80%  "

$ws.Range("D43").Value = "0.376"
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("D45").Value = "281.49"
$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").Value = "18.88"
$ws.Range("E48").Value = "  +3.40%  "

$ws.Range("D49").Value = "0.560"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("E51").Value = "  -0.85%  "

# Restore default style (clears the temporary text-number-format) while keeping values as text.
$dataRange.Style = "Normal"

